# Insert a new weekly data row for "Arveja Verde" (Vega Modelo de Temuco)
# before the existing row 60, shifting all subsequent rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(60).Insert()

$ws.Cells.Item(60, 1).Value  = 10
$ws.Cells.Item(60, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(60, 3).Value  = "La Araucanía"
$ws.Cells.Item(60, 4).Value  = 44469
$ws.Cells.Item(60, 5).Value  = 9
$ws.Cells.Item(60, 6).Value  = 100112022
$ws.Cells.Item(60, 7).Value  = "Arveja Verde"
$ws.Cells.Item(60, 8).Value  = "Sin especificar"
$ws.Cells.Item(60, 9).Value  = "Primera"
$ws.Cells.Item(60, 10).Value = 30
$ws.Cells.Item(60, 11).Value = 34000
$ws.Cells.Item(60, 12).Value = 34000
$ws.Cells.Item(60, 13).Value = 34000
$ws.Cells.Item(60, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(60, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(60, 16).Value = 1360
$ws.Cells.Item(60, 17).Value = 25
$ws.Cells.Item(60, 18).Value = "Hortaliza"
